# StructureDefinition-ror-meta-comment.xlsx edit
# 1) Metadata sheet: bump the "Date" value (row 8, column B).
# 2) Elements sheet: swap the contents of columns AK (37) and AL (38)
#    (header text + the handful of populated data cells in rows 1-6),
#    i.e. "Mapping: RIM Mapping" and
#    "Mapping: Spécification métier vers l'extension ROR MetaComment"
#    trade places.

$wb = $excel.ActiveWorkbook

# --- 1) Metadata!B8 : Date -------------------------------------------------
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Cells.Item(8, 2).Value = "2024-03-22T16:25:12+00:00"

# --- 2) Elements: swap columns AK (37) and AL (38) ------------------------
$elements = $wb.Worksheets.Item("Elements")

# Only rows 1-6 have data in this sheet (1 header row + 5 element rows).
$lastRow = 6
$colAK = 37
$colAL = 38

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elements.Cells.Item($r, $colAK)
    $alCell = $elements.Cells.Item($r, $colAL)

    $akVal = $akCell.Value2
    $alVal = $alCell.Value2

    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# The "Mapping: Spécification..." column (now AK) used to be the wide one
# (~70.07 chars) and "Mapping: RIM Mapping" (now AL) the narrow one
# (~24.98 chars) - swap the best-fit column widths to follow the content.
$elements.Columns.Item($colAK).ColumnWidth = 69.15
$elements.Columns.Item($colAL).ColumnWidth = 24.15
